$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (scraped GitHub Actions refresh).
# Column D values that look numeric get auto-converted by Excel on assignment,
# so for those we briefly force a text format, then restore the default "Normal"
# style so the cell keeps no explicit style (matching the source data, which has
# no "s" attribute on these cells) while the stored value remains a literal string.

# Row 2
$ws.Range("D2").Value = '26.658.67'
$ws.Range("E2").Value = '  +0.61%  '

# Row 3
$ws.Range("D3").Value = '1.832.01'
$ws.Range("E3").Value = '  +1.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.59%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '

# Row 7
$ws.Range("E7").Value = '  +3.46%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07147'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.62%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9314'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.47%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07649'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.13%  '

# Row 13
$ws.Range("D13").Value = '1.893.45'
$ws.Range("E13").Value = '  +4.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.259'
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = '  +0.16%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.83%  '

# Row 17
$ws.Range("E17").Value = '  +0.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008544'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.64%  '

# Row 19
$ws.Range("E19").Value = '  +0.37%  '

# Row 20
$ws.Range("D20").Value = '26.669.64'
$ws.Range("E20").Value = '  +0.54%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.19%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.022'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.90%  '

# Row 23
$ws.Range("D23").Value = '2.072.38'
$ws.Range("E23").Value = '  +0.98%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.905'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.42%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.62%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.74%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.49%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.908'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08836'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.53%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.154'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.822'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.81%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.176'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.59%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7403'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.15%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.445'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.07%  '

# Row 37
$ws.Range("E37").Value = '  +0.41%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.959'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.19%  '

# Row 39
$ws.Range("E39").Value = '  -1.13%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05153'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.85%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.921'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.62%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5065'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1499'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.111'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.91%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.008'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.47%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4662'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.26%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.72%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.77%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.575'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.18%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06026'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.72%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.11'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.15%  '
